$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.084.55"
$ws.Range("D3").Value = "2.784.01"
$ws.Range("E3").Value = "  +5.14%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'342.25"
$ws.Range("E5").Value = "  +4.46%  "
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +4.71%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +5.11%  "
$ws.Range("D10").Value = "'42.11"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("D12").Value = "'20.01"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "3.217.73"
$ws.Range("E15").Value = "  +5.03%  "
$ws.Range("D16").Value = "2.775.87"
$ws.Range("E16").Value = "  +5.23%  "
$ws.Range("D17").Value = "51.939.99"
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").Value = "'0.877"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("E19").Value = "  +9.97%  "
$ws.Range("D20").Value = "'7.04"
$ws.Range("E20").Value = "  +5.22%  "
$ws.Range("D21").Value = "'13.23"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "'277.11"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").Value = "'70.13"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").Value = "  +7.78%  "
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "'0.142"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "'34.76"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'50.22"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").Value = "'0.0820"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +3.40%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'4.97"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").Value = "'0.0384"
$ws.Range("E40").Value = "  +11.38%  "
$ws.Range("E41").Value = "  +26.87%  "
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("D44").Value = "'23.38"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "'125.97"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "2.072.02"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +5.62%  "
$ws.Range("D50").Value = "'0.907"
$ws.Range("E50").Value = "  +17.69%  "
$ws.Range("D51").Value = "'8.89"
$ws.Range("E51").Value = "  -0.71%  "
